$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete row 147 ("* Patients / Database" entry) - remaining rows shift up automatically,
# formulas (SUM range, relative refs) are auto-adjusted by Excel.
$ws.Rows("147").EntireRow.Delete()

# New shared strings must come into being in the same order the diff adds them
# (97: More realistic virtual view, 98: @IACT, 99: Debug & Refactor, 100: Bug and Crash Fixes)
# Leading "'" forces text-entry (matches the original quote-prefixed "@..." cells).
$ws.Range("E143").Value2 = [char]0x2022 + " More realistic virtual view"
$ws.Range("C137").Value = "'@IACT"
$ws.Range("B143").Value2 = "* Debug & Refactor"
$ws.Range("E145").Value2 = [char]0x2022 + " Bug and Crash Fixes"

# Remaining numeric / content updates
$ws.Range("C143").Value2 = 1
$ws.Range("C145").Value2 = 5
$ws.Range("C146").Value2 = 3
$ws.Range("E146").Clear()
$ws.Range("C147").Value2 = 28
$ws.Range("C149").Value = "'@IACT"
$ws.Range("D149").Value2 = 32

# Sheet view: scroll back to top, new selection
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("J17").Select()
